# Apply the edit described by the diff:
#  1) For every "year block" of 4 rows (A/B/C/D sub-periods), swap the
#     entire contents of the "B" sub-period row with the "C" sub-period row.
#  2) Delete columns F and G (乙烯产销率 / 乙烯销售量) entirely, header included.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ row = 2; a = '2000年A'; b = 99.59999999999999; c = $null; d = 30; e = 28.5 }
    @{ row = 3; a = '2000年C'; b = 100; c = $null; d = 3.4; e = 44.6 }
    @{ row = 4; a = '2000年B'; b = 100; c = $null; d = 0.2; e = 30.2 }
    @{ row = 5; a = '2000年D'; b = 100.1; c = $null; d = -14.3; e = 47 }
    @{ row = 6; a = '2001年A'; b = 97.2; c = $null; d = 122.2; e = 15.5 }
    @{ row = 7; a = '2001年C'; b = 99.90000000000001; c = $null; d = 34.4; e = 42.2 }
    @{ row = 8; a = '2001年B'; b = 99.7; c = $null; d = 22.1; e = 29.3 }
    @{ row = 9; a = '2001年D'; b = 99.5; c = $null; d = 76.2; e = 53.9 }
    @{ row = 10; a = '2002年A'; b = 100.9; c = $null; d = -20.6; e = 12.5 }
    @{ row = 11; a = '2002年C'; b = 100; c = $null; d = 2; e = 43.6 }
    @{ row = 12; a = '2002年B'; b = 96.8; c = $null; d = 160.5; e = 25.6 }
    @{ row = 13; a = '2002年D'; b = 100.2; c = $null; d = -10.1; e = 54.5 }
    @{ row = 14; a = '2003年A'; b = 100; c = $null; d = -0.6; e = 13.3 }
    @{ row = 15; a = '2003年C'; b = 100; c = $null; d = 0; e = 42 }
    @{ row = 16; a = '2003年B'; b = 98.2; c = $null; d = 131; e = 21.7 }
    @{ row = 17; a = '2003年D'; b = 100.1; c = $null; d = -2.4; e = 55.6 }
    @{ row = 18; a = '2004年A'; b = 100; c = $null; d = -5.1; e = 21 }
    @{ row = 19; a = '2004年C'; b = 99.90000000000001; c = $null; d = 17.9; e = 41.9 }
    @{ row = 20; a = '2004年B'; b = 99.5; c = $null; d = 41; e = 27.8 }
    @{ row = 21; a = '2004年D'; b = 100.1; c = $null; d = -20.5; e = 56.4 }
    @{ row = 22; a = '2005年A'; b = 100.4; c = $null; d = -12.9; e = 14.4 }
    @{ row = 23; a = '2005年C'; b = 99.8; c = $null; d = 26.2; e = 48.5 }
    @{ row = 24; a = '2005年B'; b = 97.40000000000001; c = $null; d = 221.4; e = 31 }
    @{ row = 25; a = '2005年D'; b = 99.8; c = $null; d = 38.2; e = 74.3 }
    @{ row = 26; a = '2006年A'; b = 99.7; c = $null; d = 10.3; e = 21.5 }
    @{ row = 27; a = '2006年C'; b = 100; c = $null; d = 2; e = 62.3 }
    @{ row = 28; a = '2006年B'; b = 99.90000000000001; c = $null; d = 4; e = 40.8 }
    @{ row = 29; a = '2006年D'; b = 99.5; c = $null; d = 88.90000000000001; e = 86.8 }
    @{ row = 30; a = '2007年A'; b = 99.40000000000001; c = $null; d = 25; e = 23.9 }
    @{ row = 31; a = '2007年C'; b = 100.1; c = $null; d = -7.1; e = 58.9 }
    @{ row = 32; a = '2007年B'; b = 99.2; c = $null; d = 5.4; e = 41.3 }
    @{ row = 33; a = '2007年D'; b = 100.3; c = $null; d = 3.6; e = 76 }
    @{ row = 34; a = '2008年A'; b = 99.5; c = $null; d = 19.3; e = 22.7 }
    @{ row = 35; a = '2008年C'; b = 99.90000000000001; c = $null; d = 3.6; e = 59.9 }
    @{ row = 36; a = '2008年B'; b = 100; c = $null; d = -4.3; e = 40.2 }
    @{ row = 37; a = '2008年D'; b = 100; c = $null; d = -5.5; e = 73.8 }
    @{ row = 38; a = '2009年A'; b = 99.8; c = $null; d = 7.1; e = 6.9 }
    @{ row = 39; a = '2009年C'; b = 99.7; c = $null; d = 42.9; e = 41.5 }
    @{ row = 40; a = '2009年B'; b = 99.8; c = $null; d = 19; e = 20 }
    @{ row = 41; a = '2009年D'; b = 99.8; c = $null; d = 46.3; e = 74.5 }
    @{ row = 42; a = '2010年A'; b = 101.3; c = $null; d = -31; e = 25.6 }
    @{ row = 43; a = '2010年C'; b = 100.5; c = $null; d = -44.6; e = 108.7 }
    @{ row = 44; a = '2010年B'; b = 101; c = $null; d = -39.3; e = 83.7 }
    @{ row = 45; a = '2010年D'; b = 100.4; c = $null; d = -55.4; e = 133.9 }
    @{ row = 46; a = '2011年A'; b = 99.3; c = $null; d = 37.5; e = 51.8 }
    @{ row = 47; a = '2011年C'; b = 99.40000000000001; c = $null; d = 133.3; e = 94.09999999999999 }
    @{ row = 48; a = '2011年B'; b = 99.7; c = $null; d = 49.1; e = 67.5 }
    @{ row = 49; a = '2011年D'; b = 99.7; c = $null; d = 95.8; e = 118.6 }
    @{ row = 50; a = '2012年A'; b = 99.8; c = $null; d = 13.9; e = 36.4 }
    @{ row = 51; a = '2012年C'; b = 100.1; c = $null; d = -11.1; e = 88.5 }
    @{ row = 52; a = '2012年B'; b = 99.90000000000001; c = $null; d = 5.6; e = 70.40000000000001 }
    @{ row = 53; a = '2012年D'; b = 100; c = $null; d = 4.2; e = 106.3 }
    @{ row = 54; a = '2016年A'; b = 99.7; c = -0.6; d = 8.4; e = 36.00314 }
    @{ row = 55; a = '2016年C'; b = 100; c = 0.1; d = 2.6; e = 119.06109 }
    @{ row = 56; a = '2016年B'; b = 99.90000000000001; c = -0.1; d = 15.8; e = 77.68309000000001 }
    @{ row = 57; a = '2016年D'; b = 100; c = -0.1; d = 1; e = 160.53691 }
    @{ row = 58; a = '2017年A'; b = 100.1; c = 0.2; d = -2.3; e = 36.44345 }
    @{ row = 59; a = '2017年C'; b = 100; c = 0.1; d = -1.5; e = 120.05571 }
    @{ row = 60; a = '2017年B'; b = 99.7; c = -0.2; d = 19.7; e = 76.49612 }
    @{ row = 61; a = '2017年D'; b = 100; c = 0.1; d = 1.6; e = 144.19058 }
    @{ row = 62; a = '2018年A'; b = 99.90000000000001; c = -0.2; d = 7.3; e = 32.02736 }
    @{ row = 63; a = '2018年C'; b = 99.90000000000001; c = -0.2; d = 20.1; e = 102.64318 }
    @{ row = 64; a = '2018年B'; b = 100; c = 0.1; d = 1.7; e = 66.27923 }
    @{ row = 65; a = '2018年D'; b = 99.90000000000001; c = -0.1; d = 19.1; e = 129.92749 }
    @{ row = 66; a = '2019年A'; b = 100.2; c = 0.3; d = -11.4; e = 44.44702 }
    @{ row = 67; a = '2019年C'; b = 99.5; c = 0.2; d = 73.3; e = 153.77427 }
    @{ row = 68; a = '2019年B'; b = 99.90000000000001; c = 0.8; d = 5; e = 100.53022 }
    @{ row = 69; a = '2019年D'; b = 100; c = 0; d = -0.7; e = 218.175 }
)

foreach ($item in $rowData) {
    $r = $item.row
    $ws.Cells.Item($r, 1).Value = $item.a
    $ws.Cells.Item($r, 2).Value = $item.b
    if ($null -ne $item.c) {
        $ws.Cells.Item($r, 3).Value = $item.c
    }
    $ws.Cells.Item($r, 4).Value = $item.d
    $ws.Cells.Item($r, 5).Value = $item.e
}

# Remove columns F:G (乙烯产销率, 乙烯销售量) for all used rows, header + data.
$ws.Range("F1:G69").Delete()

